$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 11 (2021年) below existing 2012-2020 data, years 2012-2020 occupy rows 2-10.
# Copy formatting (font/border/alignment) from A10 onto A11 before assigning the new value,
# so the new year label keeps the same bold/centered/bordered look as the other year cells.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 1828.14
$ws.Range("C11").Value = 467.98
$ws.Range("D11").Value = 87.59
$ws.Range("E11").Value = "'"
$ws.Range("F11").Value = 1191.51
$ws.Range("G11").Value = 1962.95
$ws.Range("H11").Value = 234.02
$ws.Range("I11").Value = 1425.31
$ws.Range("J11").Value = 154.13
$ws.Range("K11").Value = 33303.54
$ws.Range("L11").Value = 233.99
$ws.Range("M11").Value = 19.4
$ws.Range("N11").Value = 2.69
$ws.Range("O11").Value = 929.75
$ws.Range("P11").Value = 756.6900000000001
$ws.Range("Q11").Value = 14.1
$ws.Range("R11").Value = 82.29000000000001
$ws.Range("S11").Value = 977.74
$ws.Range("T11").Value = 31.78
$ws.Range("U11").Value = 3844.84
$ws.Range("V11").Value = "'"
$ws.Range("W11").Value = 36.14
$ws.Range("X11").Value = 147.87
$ws.Range("Y11").Value = 305.33
$ws.Range("Z11").Value = 2353.25
$ws.Range("AA11").Value = 308.39
$ws.Range("AB11").Value = 452.82
$ws.Range("AC11").Value = 42.45
$ws.Range("AD11").Value = 671.46
$ws.Range("AE11").Value = 526.26
$ws.Range("AF11").Value = 6936.22
$ws.Range("AG11").Value = 2635.6
$ws.Range("AH11").Value = 565.86
$ws.Range("AI11").Value = 365.46
$ws.Range("AJ11").Value = 88.11
$ws.Range("AK11").Value = 1165.72
$ws.Range("AL11").Value = 508.71
$ws.Range("AM11").Value = 689.0599999999999
$ws.Range("AN11").Value = 9.449999999999999
$ws.Range("AO11").Value = 545.5599999999999
$ws.Range("AP11").Value = 688.86
$ws.Range("AQ11").Value = 13.97
